$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update team member assignments (columns B:F) ---
$ws.Range("B2").Value = "EE00"
$ws.Range("C2").Value = "ME15"
$ws.Range("D2").Value = "ME27"
$ws.Range("E2").Value = "ME30"

$ws.Range("B3").Value = "EE01"
$ws.Range("C3").Value = "ME08"
$ws.Range("D3").Value = "ME31"
$ws.Range("E3").Value = "ME56"

$ws.Range("B4").Value = "EE02"
$ws.Range("C4").Value = "EE23"
$ws.Range("D4").Value = "ME29"
$ws.Range("E4").Value = "ME32"

$ws.Range("B5").Value = "EE03"
$ws.Range("C5").Value = "ME03"
$ws.Range("D5").Value = "ME33"
$ws.Range("E5").Value = "ME58"

$ws.Range("B6").Value = "EE04"
$ws.Range("C6").Value = "ME05"
$ws.Range("D6").Value = "ME34"

$ws.Range("B7").Value = "EE05"
$ws.Range("C7").Value = "ME04"
$ws.Range("D7").Value = "ME13"
$ws.Range("E7").Value = "ME35"

$ws.Range("B8").Value = "EE06"
$ws.Range("C8").Value = "ME06"
$ws.Range("D8").Value = "ME25"
$ws.Range("E8").Value = "ME36"

$ws.Range("B9").Value = "EE07"
$ws.Range("C9").Value = "ME07"
$ws.Range("D9").Value = "ME26"
$ws.Range("E9").Value = "ME37"

$ws.Range("B10").Value = "EE08"
$ws.Range("C10").Value = "ME01"
$ws.Range("D10").Value = "ME38"
$ws.Range("E10").Value = "ME57"

$ws.Range("B11").Value = "EE09"
$ws.Range("C11").Value = "ME09"
$ws.Range("D11").Value = "ME39"
$ws.Range("E11").Value = "ME53"

$ws.Range("B12").Value = "EE10"
$ws.Range("C12").Value = "ME10"
$ws.Range("D12").Value = "ME40"
$ws.Range("E12").Value = "ME62"

$ws.Range("B13").Value = "EE11"
$ws.Range("C13").Value = "ME11"
$ws.Range("D13").Value = "ME24"
$ws.Range("E13").Value = "ME41"

$ws.Range("B14").Value = "EE12"
$ws.Range("C14").Value = "ME12"
$ws.Range("D14").Value = "ME42"
$ws.Range("E14").Value = "ME55"

$ws.Range("B15").Value = "CpE01"
$ws.Range("C15").Value = "CpE02"
$ws.Range("D15").Value = "EE13"
$ws.Range("E15").Value = "ME43"

$ws.Range("B16").Value = "EE14"
$ws.Range("C16").Value = "ME14"
$ws.Range("D16").Value = "ME44"
$ws.Range("E16").Value = "ME59"

$ws.Range("B17").Value = "CpE03"
$ws.Range("C17").Value = "EE15"
$ws.Range("D17").Value = "ME45"
$ws.Range("E17").Value = "ME61"

$ws.Range("B18").Value = "EE16"
$ws.Range("C18").Value = "ME02"
$ws.Range("D18").Value = "ME16"
$ws.Range("E18").Value = "ME46"

$ws.Range("B19").Value = "EE17"
$ws.Range("C19").Value = "ME17"
$ws.Range("D19").Value = "ME23"
$ws.Range("E19").Value = "ME47"
$ws.Range("F19").Value = "ME54"

$ws.Range("B20").Value = "EE18"
$ws.Range("C20").Value = "ME00"
$ws.Range("D20").Value = "ME18"
$ws.Range("E20").Value = "ME48"

$ws.Range("B21").Value = "CpE00"
$ws.Range("C21").Value = "EE19"
$ws.Range("D21").Value = "ME19"
$ws.Range("E21").Value = "ME49"

$ws.Range("B22").Value = "EE20"
$ws.Range("C22").Value = "ME20"
$ws.Range("D22").Value = "ME50"
$ws.Range("E22").Value = "ME63"

$ws.Range("B23").Value = "EE21"
$ws.Range("C23").Value = "ME21"
$ws.Range("D23").Value = "ME28"
$ws.Range("E23").Value = "ME51"

$ws.Range("B24").Value = "EE22"
$ws.Range("C24").Value = "EE24"
$ws.Range("D24").Value = "ME22"
$ws.Range("E24").Value = "ME52"
$ws.Range("F24").Value = "ME60"

$ws.Range("B25").Value = "CE00"
$ws.Range("C25").Value = "CE03"
$ws.Range("D25").Value = "CE05"
$ws.Range("E25").Value = "CE06"

$ws.Range("B26").Value = "CE01"
$ws.Range("C26").Value = "CE02"
$ws.Range("D26").Value = "CE04"

# --- Update GPA values (columns H:L) ---
$ws.Range("H2").Value = 3.4
$ws.Range("J2").Value = 2.7
$ws.Range("K2").Value = 3.0

$ws.Range("H3").Value = 2.5
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 3.6

$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 3.9
$ws.Range("K4").Value = 2.2

$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 3.3
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.8

$ws.Range("I6").Value = 2.5
$ws.Range("J6").Value = 3.4

$ws.Range("H7").Value = 3.9
$ws.Range("J7").Value = 3.3
$ws.Range("K7").Value = 2.5

$ws.Range("H8").Value = 3.0
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 3.5
$ws.Range("K8").Value = 3.6

$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 3.7
$ws.Range("J9").Value = 3.6
$ws.Range("K9").Value = 3.7

$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 2.1
$ws.Range("K10").Value = 3.7

$ws.Range("H11").Value = 2.3
$ws.Range("I11").Value = 2.9
$ws.Range("J11").Value = 3.9
$ws.Range("K11").Value = 3.3

$ws.Range("I12").Value = 3.0
$ws.Range("J12").Value = 3.0

$ws.Range("I13").Value = 3.1
$ws.Range("J13").Value = 3.4

$ws.Range("J14").Value = 2.2
$ws.Range("K14").Value = 3.5

$ws.Range("I16").Value = 2.4
$ws.Range("K16").Value = 3.9

$ws.Range("H17").Value = 3.9
$ws.Range("K17").Value = 2.1

$ws.Range("I18").Value = 2.2
$ws.Range("J18").Value = 3.6

$ws.Range("I19").Value = 3.7
$ws.Range("J19").Value = 3.3
$ws.Range("K19").Value = 3.7
$ws.Range("L19").Value = 2.4

$ws.Range("I20").Value = 3.0
$ws.Range("J20").Value = 3.8
$ws.Range("K20").Value = 2.8

$ws.Range("H21").Value = 2.6
$ws.Range("K21").Value = 3.9

$ws.Range("I22").Value = 3.0
$ws.Range("J22").Value = 3.0
$ws.Range("K22").Value = 3.3

$ws.Range("I23").Value = 2.1
$ws.Range("J23").Value = 3.8
$ws.Range("K23").Value = 3.1

$ws.Range("H24").Value = 3.6
$ws.Range("I24").Value = 3.8
$ws.Range("J24").Value = 3.2
$ws.Range("K24").Value = 3.2
$ws.Range("L24").Value = 3.0

$ws.Range("H25").Value = 3.9

$ws.Range("H26").Value = 2.7
$ws.Range("I26").Value = 3.1

Write-Output "applied edits"